# repull data, push all data, mean calculation
# Update column F (dSF) with refreshed/repulled values for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 6
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -12
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 3
$ws.Range("F8").Value = 4
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -6
$ws.Range("F14").Value = -2
$ws.Range("F15").Value = -7
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 2
